$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 16781
$ws.Range("E2").Value = 533
$ws.Range("F2").Value = 533
$ws.Range("G2").Value = 797
$ws.Range("H2").Value = 505
$ws.Range("I2").Value = 550
$ws.Range("J2").Value = -45
$ws.Range("K2").Value = 10536
$ws.Range("L2").Value = 7189
$ws.Range("M2").Value = 3347
$ws.Range("N2").Value = 2613
$ws.Range("O2").Value = 734
$ws.Range("P2").Value = 190
$ws.Range("Q2").Value = 822
$ws.Range("R2").Value = -208
$ws.Range("S2").Value = -612
$ws.Range("T2").Value = 465
$ws.Range("U2").Value = 356
$ws.Range("V2").Value = 4018
$ws.Range("W2").Value = 3.17
$ws.Range("X2").Value = 3.01
$ws.Range("Y2").Value = 23.19
$ws.Range("Z2").Value = 4.97
$ws.Range("AA2").Value = 214.79
$ws.Range("AB2").Value = 1648.67
$ws.Range("AD2").Value = 14445
$ws.Range("AE2").Value = 7017
$ws.Range("AF2").Value = 1.58
$ws.Range("AG2").Value = 1020
$ws.Range("AH2").Value = 0.92
$ws.Range("AI2").Value = 6.9
$ws.Range("AJ2").Value = 38090950
$ws.Range("AC2").ClearContents()

# Row 3
$ws.Range("D3").Value = 18465
$ws.Range("E3").Value = 395
$ws.Range("F3").Value = 395
$ws.Range("G3").Value = 312
$ws.Range("H3").Value = 122
$ws.Range("I3").Value = 203
$ws.Range("J3").Value = -81
$ws.Range("K3").Value = 9760
$ws.Range("L3").Value = 6456
$ws.Range("M3").Value = 3304
$ws.Range("N3").Value = 2896
$ws.Range("O3").Value = 408
$ws.Range("P3").Value = 190
$ws.Range("Q3").Value = 1064
$ws.Range("R3").Value = 302
$ws.Range("S3").Value = -793
$ws.Range("T3").Value = 676
$ws.Range("U3").Value = 388
$ws.Range("V3").Value = 3395
$ws.Range("W3").Value = 2.14
$ws.Range("X3").Value = 0.66
$ws.Range("Y3").Value = 7.38
$ws.Range("Z3").Value = 1.2
$ws.Range("AA3").Value = 195.43
$ws.Range("AB3").Value = 1372.32
$ws.Range("AD3").Value = 5334
$ws.Range("AE3").Value = 7776
$ws.Range("AF3").Value = 2.39
$ws.Range("AG3").Value = 1020
$ws.Range("AH3").Value = 0.55
$ws.Range("AI3").Value = 18.69
$ws.Range("AJ3").Value = 38090950
$ws.Range("AC3").ClearContents()

# Row 4
$ws.Range("D4").Value = 20307
$ws.Range("E4").Value = 379
$ws.Range("F4").Value = 379
$ws.Range("G4").Value = 289
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 186
$ws.Range("J4").Value = -86
$ws.Range("K4").Value = 10614
$ws.Range("L4").Value = 6951
$ws.Range("M4").Value = 3662
$ws.Range("N4").Value = 2602
$ws.Range("O4").Value = 1061
$ws.Range("P4").Value = 190
$ws.Range("Q4").Value = 579
$ws.Range("R4").Value = -1487
$ws.Range("S4").Value = 308
$ws.Range("T4").Value = 694
$ws.Range("U4").Value = -115
$ws.Range("V4").Value = 3467
$ws.Range("W4").Value = 1.87
$ws.Range("X4").Value = 0.49
$ws.Range("Y4").Value = 6.77
$ws.Range("Z4").Value = 0.98
$ws.Range("AA4").Value = 189.81
$ws.Range("AB4").Value = 1249.51
$ws.Range("AD4").Value = 4882
$ws.Range("AE4").Value = 6986
$ws.Range("AF4").Value = 2.01
$ws.Range("AG4").Value = 1020
$ws.Range("AH4").Value = 0.73
$ws.Range("AI4").Value = 20.43
$ws.Range("AJ4").Value = 38090950
$ws.Range("AC4").ClearContents()

# Row 5
$ws.Range("D5").Value = 21957
$ws.Range("E5").Value = 528
$ws.Range("F5").Value = 528
$ws.Range("G5").Value = 517
$ws.Range("H5").Value = 304
$ws.Range("I5").Value = 385
$ws.Range("J5").Value = -81
$ws.Range("K5").Value = 11030
$ws.Range("L5").Value = 7272
$ws.Range("M5").Value = 3758
$ws.Range("N5").Value = 2678
$ws.Range("O5").Value = 1080
$ws.Range("P5").Value = 190
$ws.Range("Q5").Value = 626
$ws.Range("R5").Value = -950
$ws.Range("S5").Value = 192
$ws.Range("T5").Value = 1035
$ws.Range("U5").Value = -408
$ws.Range("V5").Value = 3691
$ws.Range("W5").Value = 2.4
$ws.Range("X5").Value = 1.39
$ws.Range("Y5").Value = 14.59
$ws.Range("Z5").Value = 2.81
$ws.Range("AA5").Value = 193.49
$ws.Range("AB5").Value = 1344.68
$ws.Range("AC5").Value = 1011
$ws.Range("AD5").Value = 17.01
$ws.Range("AE5").Value = 7192
$ws.Range("AF5").Value = 2.39
$ws.Range("AG5").Value = 1153
$ws.Range("AH5").Value = 0.67
$ws.Range("AI5").Value = 11.15
$ws.Range("AJ5").Value = 38090950

# Row 6
$ws.Range("D6").Value = 22720
$ws.Range("E6").Value = 402
$ws.Range("F6").Value = 402
$ws.Range("G6").Value = 276
$ws.Range("H6").Value = 126
$ws.Range("I6").Value = 234
$ws.Range("K6").Value = 12146
$ws.Range("L6").Value = 7748
$ws.Range("M6").Value = 4398
$ws.Range("N6").Value = 3311
$ws.Range("P6").Value = 211
$ws.Range("Q6").Value = 593
$ws.Range("R6").Value = -1365
$ws.Range("S6").Value = 1071
$ws.Range("T6").Value = 914
$ws.Range("U6").Value = -321
$ws.Range("V6").Value = 4369
$ws.Range("W6").Value = 1.77
$ws.Range("X6").Value = 0.5600000000000001
$ws.Range("Y6").Value = 7.82
$ws.Range("Z6").Value = 1.09
$ws.Range("AA6").Value = 176.18
$ws.Range("AB6").Value = 1481.22
$ws.Range("AC6").Value = 568
$ws.Range("AD6").Value = 14.02
$ws.Range("AE6").Value = 8021
$ws.Range("AF6").Value = 0.99
$ws.Range("AG6").Value = 1020
$ws.Range("AH6").Value = 1.28
$ws.Range("AI6").Value = 22.2
$ws.Range("AJ6").Value = 38090950

# Row 7
$ws.Range("D7").Value = 23661
$ws.Range("E7").Value = 322
$ws.Range("G7").Value = 187
$ws.Range("H7").Value = 84
$ws.Range("I7").Value = 159
$ws.Range("K7").Value = 14582
$ws.Range("L7").Value = 10066
$ws.Range("M7").Value = 4516
$ws.Range("N7").Value = 3598
$ws.Range("P7").Value = 210
$ws.Range("Q7").Value = 891
$ws.Range("R7").Value = -1067
$ws.Range("S7").Value = -7
$ws.Range("T7").Value = 930
$ws.Range("U7").Value = -148
$ws.Range("W7").Value = 1.36
$ws.Range("X7").Value = 0.36
$ws.Range("Y7").Value = 4.6
$ws.Range("Z7").Value = 0.63
$ws.Range("AA7").Value = 222.89
$ws.Range("AC7").Value = 377
$ws.Range("AD7").Value = 28.35
$ws.Range("AE7").Value = 8718
$ws.Range("AF7").Value = 1.23
$ws.Range("AG7").Value = 81
$ws.Range("AH7").Value = 0.75
$ws.Range("AI7").Value = 19.32

# Row 8
$ws.Range("D8").Value = 24710
$ws.Range("E8").Value = 460
$ws.Range("G8").Value = 330
$ws.Range("H8").Value = 209
$ws.Range("I8").Value = 237
$ws.Range("K8").Value = 14843
$ws.Range("L8").Value = 10121
$ws.Range("M8").Value = 4722
$ws.Range("N8").Value = 3799
$ws.Range("P8").Value = 210
$ws.Range("Q8").Value = 1086
$ws.Range("R8").Value = -1022
$ws.Range("S8").Value = -133
$ws.Range("T8").Value = 865
$ws.Range("U8").Value = 140
$ws.Range("W8").Value = 1.86
$ws.Range("X8").Value = 0.84
$ws.Range("Y8").Value = 6.4
$ws.Range("Z8").Value = 1.42
$ws.Range("AA8").Value = 214.34
$ws.Range("AC8").Value = 562
$ws.Range("AD8").Value = 18.78
$ws.Range("AE8").Value = 9204
$ws.Range("AF8").Value = 1.15
$ws.Range("AG8").Value = 91
$ws.Range("AH8").Value = 0.86
$ws.Range("AI8").Value = 14.59

# Row 9
$ws.Range("D9").Value = 26272
$ws.Range("E9").Value = 563
$ws.Range("G9").Value = 448
$ws.Range("H9").Value = 270
$ws.Range("I9").Value = 309
$ws.Range("K9").Value = 15131
$ws.Range("L9").Value = 10134
$ws.Range("M9").Value = 4996
$ws.Range("N9").Value = 4075
$ws.Range("P9").Value = 210
$ws.Range("Q9").Value = 1152
$ws.Range("R9").Value = -1037
$ws.Range("S9").Value = -226
$ws.Range("T9").Value = 919
$ws.Range("U9").Value = 260
$ws.Range("W9").Value = 2.14
$ws.Range("X9").Value = 1.03
$ws.Range("Y9").Value = 7.84
$ws.Range("Z9").Value = 1.8
$ws.Range("AA9").Value = 202.84
$ws.Range("AC9").Value = 733
$ws.Range("AD9").Value = 14.4
$ws.Range("AE9").Value = 9873
$ws.Range("AF9").Value = 1.07
$ws.Range("AG9").Value = 104
$ws.Range("AH9").Value = 0.99
$ws.Range("AI9").Value = 12.83

